# Adds two new client/property/financial/status rows (033/2025 and 034/2025)
# to the bottom of each of the four data sheets, mirroring the existing
# row-32/33 pattern (values + styles + hyperlinks), per the commit:
# "feat: adiciona a interface grafica e algumas melhorias extras"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Dados dos Clientes" -> add rows 34 and 35
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Dados dos Clientes")

# Row 34 - Daniely Evellin
$ws1.Range("A33:E33").Copy($ws1.Range("A34:E34"))
$ws1.Range("A34").Value = "033/2025"
$ws1.Range("B34").Value = "Daniely Evellin"
$ws1.Range("C34").Value = "456.789.123-55"
$ws1.Range("D34").Value = "dani.evellin@email.com"
$ws1.Hyperlinks.Add($ws1.Range("D34"), "mailto:dani.evellin@email.com")
$ws1.Range("D34").Style = $ws1.Range("D33").Style

# Row 35 - Erika Polina (document/CPF entered before the name, matching
# the original authoring order captured in the shared-strings table)
$ws1.Range("A33:E33").Copy($ws1.Range("A35:E35"))
$ws1.Range("A35").Value = "034/2025"
$ws1.Range("C35").Value = "123.456.789-55"
$ws1.Range("B35").Value = "Erika Polina"
$ws1.Range("D35").Value = "erika.poliana@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("D35"), "mailto:erika.poliana@gmail.com")
$ws1.Range("D35").Style = $ws1.Range("D33").Style

[void]$ws1.Range("F39").Select()

# ---------------------------------------------------------------------
# Sheet 2: "Dados do ImóvelServiço Contrata" -> add rows 34 and 35
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Dados do ImóvelServiço Contrata")

$ws2.Range("A33:D33").Copy($ws2.Range("A34:D34"))
$ws2.Range("A33:D33").Copy($ws2.Range("A35:D35"))

[void]$ws2.Range("A35:D35").Select()

# ---------------------------------------------------------------------
# Sheet 3: "Dados Financeiros e Valores" -> add rows 34 and 35
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Dados Financeiros e Valores")

$ws3.Range("A33:D33").Copy($ws3.Range("A34:D34"))
$ws3.Range("A33:D33").Copy($ws3.Range("A35:D35"))

[void]$ws3.Range("A35:D35").Select()

# ---------------------------------------------------------------------
# Sheet 4: "StatusErro" -> add rows 34 and 35
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("StatusErro")

$ws4.Range("A33:B33").Copy($ws4.Range("A34:B34"))
$ws4.Range("A33:B33").Copy($ws4.Range("A35:B35"))

[void]$ws4.Range("A35:B35").Select()

[void]$ws1.Activate()
